$wb = $excel.ActiveWorkbook

# --- Sheet "FBS" (sheet1) ---
$ws1 = $wb.Worksheets.Item("FBS")

$ws1.Range("O4").Value = 84.44000000000001
$ws1.Range("Q15").Value = "SE"
$ws1.Range("R19").Value = 0.9000000000000001
$ws1.Range("O20").Value = 48.92000000000001
$ws1.Range("Q24").Value = "NE"
$ws1.Range("O27").Value = 55.84999999999999
$ws1.Range("O29").Value = 65.11999999999999
$ws1.Range("Q32").Value = "NNE"
$ws1.Range("Q33").Value = "SE"
$ws1.Range("O34").Value = 66.61999999999999
$ws1.Range("Q35").Value = "SE"
$ws1.Range("Q36").Value = "S"
$ws1.Range("O47").Value = 65.11999999999999

# --- Sheet "Other" (sheet2) ---
$ws2 = $wb.Worksheets.Item("Other")

$ws2.Range("Q22").Value = 67.75999999999999
$ws2.Range("Q23").Value = 46.09999999999999
$ws2.Range("S23").Value = "SE"
$ws2.Range("Q25").Value = 48.74000000000001
$ws2.Range("Q37").Value = 67.09999999999999
$ws2.Range("Q41").Value = 79.81999999999999
$ws2.Range("Q48").Value = 55.09999999999999

# --- Update Timestamp column (shared string reused across every data row, column AK on FBS sheet) ---
$oldTimestamp = "2024-11-07T05:12:44.964176"
$newTimestamp = "2024-11-07T05:15:54.442111"

$lastRow = $ws1.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws1.Cells.Item($r, 37) # column AK
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
